$d = $word.ActiveDocument

# 1. "...new ways fallen faith..." -> "...new ways and fallen faith..."
$d.Content.Find.Execute("new ways fallen faith", $true, $false, $false, $false, $false, `
    $true, 1, $false, "new ways and fallen faith", 2) | Out-Null

# 2. Fix typo "leathal" -> "lethal"
$d.Content.Find.Execute("dropping monstrous, and leathal,", $true, $false, $false, $false, $false, `
    $true, 1, $false, "dropping monstrous, and lethal,", 2) | Out-Null

# 3. Remove double space before "temperature is always cold"
$d.Content.Find.Execute("the  temperature is always cold", $true, $false, $false, $false, $false, `
    $true, 1, $false, "the temperature is always cold", 2) | Out-Null

# 4. Capitalize "ice mountains" -> "Ice Mountains"
$d.Content.Find.Execute("marked with ice mountains and", $true, $false, $false, $false, $false, `
    $true, 1, $false, "marked with Ice Mountains and", 2) | Out-Null

# 5. Fix "it's edges" -> "its edges" (possessive, no apostrophe)
$d.Content.Find.Execute([string]::Concat("either of it", [char]0x2019, "s edges"), $true, $false, $false, $false, $false, `
    $true, 1, $false, "either of its edges", 2) | Out-Null

# 6. Fix typo "develope" -> "develop", add commas: "mining and metal weaponry and tools" -> "mining, metal weaponry, and tools"
$d.Content.Find.Execute("develope around mining and metal weaponry and tools.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "develop around mining, metal weaponry, and tools.", 2) | Out-Null

# 7. Fix "froze" -> "frozen"
$d.Content.Find.Execute("the ground is constantly froze and there is no light", $true, $false, $false, $false, $false, `
    $true, 1, $false, "the ground is constantly frozen and there is no light", 2) | Out-Null

# 8. Fix "centered around" -> "centered on"
$d.Content.Find.Execute("past times are centered around social engagements", $true, $false, $false, $false, $false, `
    $true, 1, $false, "past times are centered on social engagements", 2) | Out-Null

# 9. Remove the stray "Change you Git!" paragraph together with one of the two blank
#    paragraphs that preceded it, leaving a single blank paragraph in their place.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Change you Git!") {
        $gitIndex = $i
        break
    }
}
$pGit = $d.Paragraphs.Item($gitIndex)
$rngGit = $d.Range($pGit.Range.Start, $pGit.Range.End)
$rngGit.Delete() | Out-Null

$pBlank = $d.Paragraphs.Item($gitIndex - 1)
$rngBlank = $d.Range($pBlank.Range.Start, $pBlank.Range.End)
$rngBlank.Delete() | Out-Null
